$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.549.43"
$ws.Range("E2").Value = "  +1.43%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.707.30"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "616.11"
$ws.Range("E5").Value = "  +6.47%  "

# Row 6 - Solana
$ws.Range("D6").Value = "196.23"
$ws.Range("E6").Value = "  +14.85%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.636"
$ws.Range("E7").Value = "  +2.23%  "

# Row 8 - USDC
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +3.68%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "61.08"
$ws.Range("E10").Value = "  +19.09%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.28%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  -0.90%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +0.25%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.294.88"
$ws.Range("E14").Value = "  +0.99%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.705.61"
$ws.Range("E15").Value = "  +0.04%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "19.49"
$ws.Range("E16").Value = "  +0.96%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +0.95%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +3.11%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +0.54%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "68.404.94"
$ws.Range("E20").Value = "  +1.48%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "409.16"
$ws.Range("E21").Value = "  +1.34%  "

# Row 22 - PancakeSwap
$ws.Range("D22").Value = "4.67"
$ws.Range("E22").Value = "  +4.35%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "89.97"
$ws.Range("E23").Value = "  +2.99%  "

# Row 24 - now RenderToken (was ImmutableX)
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "11.56"
$ws.Range("E24").Value = "  +9.29%  "

# Row 25 - now ImmutableX (was RenderToken)
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "3.08"
$ws.Range("E25").Value = "  +1.69%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "13.07"
$ws.Range("E26").Value = "  +2.90%  "

# Row 27 - LEO
$ws.Range("D27").Value = "6.04"
$ws.Range("E27").Value = "  +1.02%  "

# Row 28 - Toncoin
$ws.Range("D28").Value = "3.79"
$ws.Range("E28").Value = "  +1.59%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  +2.79%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "32.76"
$ws.Range("E30").Value = "  +0.87%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.68"
$ws.Range("E31").Value = "  +2.97%  "

# Row 32 - now InjectiveProtocol (was Cosmos)
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "47.82"
$ws.Range("E32").Value = "  +11.47%  "

# Row 33 - now Cosmos (was InjectiveProtocol)
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "12.69"
$ws.Range("E33").Value = "  +2.25%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +5.98%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "636.32"
$ws.Range("E35").Value = "  +7.37%  "

# Row 36 - OKB
$ws.Range("D36").Value = "67.57"
$ws.Range("E36").Value = "  +4.18%  "

# Row 37 - TheGraph
$ws.Range("D37").Value = "0.413"
$ws.Range("E37").Value = "  +4.93%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0819"
$ws.Range("E38").Value = "  -6.83%  "

# Row 39 - Dai
$ws.Range("E39").Value = "  +0.19%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.14%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +4.58%  "

# Row 42 - ThetaToken
$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  +2.40%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "0.0445"
$ws.Range("E43").Value = "  +2.31%  "

# Row 44 - Fetch.AI
$ws.Range("D44").Value = "2.62"
$ws.Range("E44").Value = "  -0.83%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.927.24"
$ws.Range("E45").Value = "  +5.38%  "

# Row 46 - THORChain
$ws.Range("D46").Value = "9.43"
$ws.Range("E46").Value = "  +3.12%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  +4.95%  "

# Row 48 - Monero
$ws.Range("D48").Value = "146.73"
$ws.Range("E48").Value = "  +2.70%  "

# Row 50 - ApeXProtocol
$ws.Range("D50").Value = "3.07"
$ws.Range("E50").Value = "  -2.41%  "

# Row 51 - dogwifhat
$ws.Range("D51").Value = "2.58"
$ws.Range("E51").Value = "  -11.39%  "
